$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# shared string "2017-02-22 08:41:59" -> "2017-02-22 08:44:07"
# used by Overview!G2:G3 (Latest HO Xliff Generate Date)
# and de-de!H2:H3 (Correspond Handoff Datetime)
$wsOverview.Range("G2").Value = "2017-02-22 08:44:07"
$wsOverview.Range("G3").Value = "2017-02-22 08:44:07"
$wsDeDe.Range("H2").Value = "2017-02-22 08:44:07"
$wsDeDe.Range("H3").Value = "2017-02-22 08:44:07"

# shared string "ht" -> "mt"
# used by zh-cn!E2:E3 and de-de!E2:E3 (Priority)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# shared string "2017-02-22 08:41:41" -> "2017-02-22 08:43:50"
# used by zh-cn!H2:H3 (Correspond Handoff Datetime)
$wsZhCn.Range("H2").Value = "2017-02-22 08:43:50"
$wsZhCn.Range("H3").Value = "2017-02-22 08:43:50"

# shared string "2017-02-22 08:42:40" -> "2017-02-22 08:44:49"
# used by zh-cn!L2:L3 (Correspond Handback DateTime)
$wsZhCn.Range("L2").Value = "2017-02-22 08:44:49"
$wsZhCn.Range("L3").Value = "2017-02-22 08:44:49"

# shared string "2017-02-22 08:43:03" -> "2017-02-22 08:45:13"
# used by de-de!L2:L3 (Correspond Handback DateTime)
$wsDeDe.Range("L2").Value = "2017-02-22 08:45:13"
$wsDeDe.Range("L3").Value = "2017-02-22 08:45:13"
